$wb = $excel.ActiveWorkbook

# --- Content edit: replace "old dwarf mine" treasure-room quest entry with the
# --- new "Snakesign tavern / drunken NPC" quest entry, and add a second row
# --- describing the tavern itself. This is the real authoring change behind
# --- the commit ("Improved Karl cave quest").

$ws3 = $wb.Worksheets.Item("GlobalVars")
$ws3.Cells.Item(9, 1).Value = "226: You visited the tavern in Snakesign"

$wsMap = $wb.Worksheets.Item("MapChanges")
$wsMap.Cells.Item(10, 2).Value = "Tavern of the goddess"
$wsMap.Cells.Item(9, 2).Value = "Snakesign"
$wsMap.Cells.Item(9, 3).Value = "Added drunken NPC"
$wsMap.Cells.Item(10, 3).Value = "Global var 226 is now set when you enter the tavern, added NPC who talks about Karl"
$wsMap.Cells.Item(9, 1).Value = 420
$wsMap.Cells.Item(10, 1).Value = 421
$wsMap.Rows.Item(10).RowHeight = 30

# --- View-state changes: the author ended the session with MapChanges active
# --- (having navigated there from Todo), with specific cell selections left
# --- on GlobalVars and MapChanges.

$ws3.Activate()
$ws3.Range("A10").Select()

$wsMap.Activate()
$wsMap.Range("A11").Select()
